$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44335
$ws.Range("M2").Value = 80
$ws.Range("N2").Value = 14000
$ws.Range("O2").Value = 14000
$ws.Range("P2").Value = 14000
$ws.Range("R2").Value = 'Región Metropolitana'
$ws.Range("S2").Value = 2000
$ws.Range("D3").Value = 44342
$ws.Range("M3").Value = 50
$ws.Range("D4").Value = 44664
$ws.Range("M4").Value = 80
$ws.Range("N4").Value = 14000
$ws.Range("O4").Value = 14000
$ws.Range("P4").Value = 14000
$ws.Range("S4").Value = 2000
$ws.Range("D5").Value = 44664
$ws.Range("M5").Value = 50
$ws.Range("D6").Value = 44344
$ws.Range("L6").Value = 'Segunda'
$ws.Range("M6").Value = 50
$ws.Range("N6").Value = 12000
$ws.Range("O6").Value = 12000
$ws.Range("P6").Value = 12000
$ws.Range("S6").Value = 1714
$ws.Range("D7").Value = 44659
$ws.Range("M7").Value = 50
$ws.Range("D8").Value = 44659
$ws.Range("M8").Value = 20
$ws.Range("D9").Value = 44300
$ws.Range("M9").Value = 150
$ws.Range("N9").Value = 12000
$ws.Range("P9").Value = 12500
$ws.Range("R9").Value = 'Provincia de Santiago'
$ws.Range("S9").Value = 1786
$ws.Range("D10").Value = 44679
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 150
$ws.Range("N10").Value = 12000
$ws.Range("O10").Value = 12000
$ws.Range("P10").Value = 12000
$ws.Range("S10").Value = 1714
$ws.Range("D11").Value = 44623
$ws.Range("L11").Value = 'Segunda'
$ws.Range("M11").Value = 30
$ws.Range("N11").Value = 16000
$ws.Range("O11").Value = 16000
$ws.Range("P11").Value = 16000
$ws.Range("S11").Value = 2286
$ws.Range("D12").Value = 44685
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 100
$ws.Range("N12").Value = 15000
$ws.Range("O12").Value = 15000
$ws.Range("P12").Value = 15000
$ws.Range("S12").Value = 2143
$ws.Range("D13").Value = 44685
$ws.Range("L13").Value = 'Segunda'
$ws.Range("M13").Value = 70
$ws.Range("D14").Value = 44687
$ws.Range("L14").Value = 'Primera'
$ws.Range("M14").Value = 100
$ws.Range("N14").Value = 15000
$ws.Range("O14").Value = 15000
$ws.Range("P14").Value = 15000
$ws.Range("S14").Value = 2143
$ws.Range("D15").Value = 44687
$ws.Range("M15").Value = 75
$ws.Range("N15").Value = 12000
$ws.Range("O15").Value = 12000
$ws.Range("P15").Value = 12000
$ws.Range("S15").Value = 1714
$ws.Range("D16").Value = 44321
$ws.Range("M16").Value = 140
$ws.Range("N16").Value = 11000
$ws.Range("O16").Value = 12000
$ws.Range("P16").Value = 11500
$ws.Range("S16").Value = 1643
$ws.Range("D17").Value = 44321
$ws.Range("M17").Value = 80
$ws.Range("N17").Value = 8000
$ws.Range("O17").Value = 8000
$ws.Range("P17").Value = 8000
$ws.Range("S17").Value = 1143
$ws.Range("D18").Value = 44312
$ws.Range("M18").Value = 50
$ws.Range("D19").Value = 44312
$ws.Range("M19").Value = 20
$ws.Range("D20").Value = 44307
$ws.Range("M20").Value = 70
$ws.Range("N20").Value = 14000
$ws.Range("O20").Value = 14000
$ws.Range("P20").Value = 14000
$ws.Range("R20").Value = 'Región Metropolitana'
$ws.Range("S20").Value = 2000
$ws.Range("D21").Value = 44307
$ws.Range("L21").Value = 'Segunda'
$ws.Range("M21").Value = 50
$ws.Range("N21").Value = 10000
$ws.Range("O21").Value = 10000
$ws.Range("P21").Value = 10000
$ws.Range("S21").Value = 1429
$ws.Range("D22").Value = 44690
$ws.Range("N22").Value = 12000
$ws.Range("O22").Value = 12000
$ws.Range("P22").Value = 12000
$ws.Range("S22").Value = 1714
$ws.Range("D23").Value = 44641
$ws.Range("M23").Value = 100
$ws.Range("N23").Value = 13000
$ws.Range("O23").Value = 13000
$ws.Range("P23").Value = 13000
$ws.Range("S23").Value = 1857
$ws.Range("D24").Value = 44306
$ws.Range("L24").Value = 'Primera'
$ws.Range("N24").Value = 12000
$ws.Range("O24").Value = 12000
$ws.Range("P24").Value = 12000
$ws.Range("S24").Value = 1714
$ws.Range("D25").Value = 44306
$ws.Range("L25").Value = 'Segunda'
$ws.Range("M25").Value = 40
$ws.Range("N25").Value = 9000
$ws.Range("O25").Value = 9000
$ws.Range("P25").Value = 9000
$ws.Range("S25").Value = 1286
$ws.Range("D26").Value = 44657
$ws.Range("M26").Value = 100
$ws.Range("N26").Value = 13000
$ws.Range("P26").Value = 13000
$ws.Range("S26").Value = 1857
$ws.Range("D27").Value = 44314
$ws.Range("L27").Value = 'Primera'
$ws.Range("M27").Value = 20
$ws.Range("N27").Value = 13000
$ws.Range("O27").Value = 13000
$ws.Range("P27").Value = 13000
$ws.Range("S27").Value = 1857
$ws.Range("D28").Value = 44314
$ws.Range("L28").Value = 'Segunda'
$ws.Range("M28").Value = 45
$ws.Range("N28").Value = 11000
$ws.Range("O28").Value = 11000
$ws.Range("P28").Value = 11000
$ws.Range("S28").Value = 1571
$ws.Range("D29").Value = 44315
$ws.Range("L29").Value = 'Especial'
$ws.Range("N29").Value = 14000
$ws.Range("O29").Value = 14000
$ws.Range("P29").Value = 14000
$ws.Range("S29").Value = 2000
$ws.Range("D30").Value = 44315
$ws.Range("M30").Value = 80
$ws.Range("N30").Value = 12000
$ws.Range("O30").Value = 13000
$ws.Range("P30").Value = 12500
$ws.Range("S30").Value = 1786
$ws.Range("D31").Value = 44315
$ws.Range("L31").Value = 'Segunda'
$ws.Range("N31").Value = 10000
$ws.Range("O31").Value = 11000
$ws.Range("P31").Value = 10500
$ws.Range("S31").Value = 1500
$ws.Range("D32").Value = 44302
$ws.Range("M32").Value = 340
$ws.Range("N32").Value = 12000
$ws.Range("O32").Value = 13000
$ws.Range("P32").Value = 12500
$ws.Range("R32").Value = 'Provincia de Santiago'
$ws.Range("S32").Value = 1786
$ws.Range("D33").Value = 44349
$ws.Range("L33").Value = 'Segunda'
$ws.Range("M33").Value = 70
$ws.Range("N33").Value = 12000
$ws.Range("O33").Value = 12000
$ws.Range("P33").Value = 12000
$ws.Range("S33").Value = 1714
$ws.Range("D34").Value = 44316
$ws.Range("L34").Value = 'Primera'
$ws.Range("M34").Value = 40
$ws.Range("N34").Value = 13000
$ws.Range("O34").Value = 13000
$ws.Range("P34").Value = 13000
$ws.Range("S34").Value = 1857
$ws.Range("D35").Value = 44316
$ws.Range("L35").Value = 'Segunda'
$ws.Range("M35").Value = 50
$ws.Range("N35").Value = 11000
$ws.Range("O35").Value = 11000
$ws.Range("P35").Value = 11000
$ws.Range("S35").Value = 1571
$ws.Range("D36").Value = 44322
$ws.Range("L36").Value = 'Primera'
$ws.Range("M36").Value = 100
$ws.Range("D37").Value = 44643
$ws.Range("M37").Value = 100
$ws.Range("N37").Value = 15000
$ws.Range("O37").Value = 15000
$ws.Range("P37").Value = 15000
$ws.Range("S37").Value = 2143
$ws.Range("D38").Value = 44644
$ws.Range("L38").Value = 'Primera'
$ws.Range("M38").Value = 85
$ws.Range("N38").Value = 14000
$ws.Range("O38").Value = 14000
$ws.Range("P38").Value = 14000
$ws.Range("S38").Value = 2000
